$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Site"
$ws.Range("B1").Value = "Login"
$ws.Range("C1").Value = "Passwords"

$ws.Range("A2").Value = "www.d"
$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "Z.=QDH6AfP}E"
